$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Add()
$ws2.Range("A1").Value = "Col1"
$ws2.Range("B1").Value = "Col2"
$ws2.Range("A2").Value = "x"
$ws2.Range("B2").Value = "y"
$lo = $ws2.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws2.Range("A1:B2"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
Write-Host "Added:" $lo.Name
# Apply styling AFTER table creation
$ws2.Range("A1:B1").Font.Bold = $true
$ws2.Range("A1:B1").Interior.Color = RGB(217,217,217)
$ws2.Range("A1:B1").Borders.LineStyle = 1
